# Add a new derived-variable row ("Rx15" / ac_apa_baseline / Baseline
# anticoagulation, aspirin, or APA) into the Treatments section of the
# table, right above the existing "Rx1" (row 57) entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new blank row at sheet row 57; this shifts every row below
# (the rest of the table) down by one, including the last table row.
$ws.Rows.Item(57).Insert()

# Populate the new row with the new variable's data.
$ws.Range("A57").Value = "Rx15"
$ws.Range("B57").Value = "ac_apa_baseline"
$ws.Range("C57").Value = "Treatments"
$ws.Range("D57").Value = "Baseline anticoagulation, aspirin, or APA"

# The table definition needs to grow by one row to re-absorb the row
# that got pushed past its previous bottom edge.
$lo.Resize($ws.Range("A1:E85"))

# Reflect where the author ended up after making the edit.
$ws.Range("D57").Select()
